$wb = $excel.ActiveWorkbook

# --- Step1_Data ---
$ws = $wb.Worksheets.Item("Step1_Data")
$ws.Range("B2").Value = 0
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = 0.1608144049558392
$ws.Range("E2").Value = 0
$ws.Range("F2").Value = 0.3703727364680821
$ws.Range("G2").Value = 0
$ws.Range("H2").Value = 0.009818955953831142
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 0.04293372513123797
$ws.Range("M2").Value = 0.006775480508511659
$ws.Range("N2").Value = 0.112460354537128
$ws.Range("O2").Value = 0
$ws.Range("P2").Value = 0.01102242327698758
$ws.Range("Q2").Value = 0
$ws.Range("R2").Value = 0.06299278693095725
$ws.Range("S2").Value = 0.05594106352041486
$ws.Range("T2").Value = 0.007073425420849847
$ws.Range("U2").Value = 0
$ws.Range("V2").Value = 0.0438218068010067
$ws.Range("W2").Value = 0
$ws.Range("X2").Value = 0.00916606451412016
$ws.Range("Y2").Value = 0.07559312628516208
$ws.Range("Z2").Value = 0.03121364569587169
$ws.Range("AA2").Value = 0
$ws.Range("AB2").Value = 0
$ws.Range("AC2").Value = 0
$ws.Range("AD2").Value = 0
$ws.Range("AE2").Value = 0
$ws.Range("AF2").Value = 0
$ws.Range("AG2").Value = 0
$ws.Range("AH2").Value = 0
$ws.Range("AI2").Value = 0
$ws.Range("B3").Value = 0
$ws.Range("C3").Value = 0
$ws.Range("D3").Value = 0.09328758731335501
$ws.Range("E3").Value = 0
$ws.Range("F3").Value = 0.3467044346567913
$ws.Range("G3").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 0.0431195346405815
$ws.Range("M3").Value = 0
$ws.Range("N3").Value = 0.112474517035731
$ws.Range("O3").Value = 0.01918185967953922
$ws.Range("P3").Value = 0
$ws.Range("Q3").Value = 0
$ws.Range("R3").Value = 0.01573615402619858
$ws.Range("S3").Value = 0.08256452044515905
$ws.Range("T3").Value = 0
$ws.Range("U3").Value = 0.0270611468869987
$ws.Range("V3").Value = 0.03409372586944558
$ws.Range("W3").Value = 0
$ws.Range("X3").Value = 0.06412528834009124
$ws.Range("Y3").Value = 0.02220976361470356
$ws.Range("Z3").Value = 0.1249589125895823
$ws.Range("AA3").Value = 0
$ws.Range("AB3").Value = 0.011453950561365
$ws.Range("AC3").Value = 0
$ws.Range("AD3").Value = 0
$ws.Range("AE3").Value = 0.003028604340457946
$ws.Range("AF3").Value = 0
$ws.Range("AG3").Value = 0
$ws.Range("AH3").Value = 0
$ws.Range("AI3").Value = 0
$ws.Range("B4").Value = 0
$ws.Range("C4").Value = 0
$ws.Range("D4").Value = 0.2229819650632389
$ws.Range("E4").Value = 0.054140272480428
$ws.Range("F4").Value = 0.2343514060212586
$ws.Range("G4").Value = 0
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 0.04897349854806408
$ws.Range("M4").Value = 0.01680768392020967
$ws.Range("N4").Value = 0.0961200363568406
$ws.Range("O4").Value = 0
$ws.Range("P4").Value = 0.003878859088935221
$ws.Range("Q4").Value = 0
$ws.Range("R4").Value = 0.06247691805989709
$ws.Range("S4").Value = 0.03349370861501081
$ws.Range("T4").Value = 0.01977598033446267
$ws.Range("U4").Value = 0
$ws.Range("V4").Value = 0.05774155128235835
$ws.Range("W4").Value = 0.006626287458440335
$ws.Range("X4").Value = 0.008819952274019221
$ws.Range("Y4").Value = 0.1065380672768528
$ws.Range("Z4").Value = 0.0211332044502668
$ws.Range("AA4").Value = 0.006140608769716525
$ws.Range("AB4").Value = 0
$ws.Range("AC4").Value = 0
$ws.Range("AD4").Value = 0
$ws.Range("AE4").Value = 0
$ws.Range("AF4").Value = 0
$ws.Range("AG4").Value = 0
$ws.Range("AH4").Value = 0
$ws.Range("AI4").Value = 0
$ws.Range("B5").Value = 0
$ws.Range("C5").Value = 0
$ws.Range("D5").Value = 0.1235929690994866
$ws.Range("E5").Value = 0
$ws.Range("F5").Value = 0.3523841177464215
$ws.Range("G5").Value = 0
$ws.Range("H5").Value = 0.0109819547268525
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 0.03530998264615263
$ws.Range("M5").Value = 0
$ws.Range("N5").Value = 0.1316385683445145
$ws.Range("O5").Value = 0
$ws.Range("P5").Value = 0
$ws.Range("Q5").Value = 0.006045082988199479
$ws.Range("R5").Value = 0.02590911826532901
$ws.Range("S5").Value = 0.08817835133756781
$ws.Range("T5").Value = 0
$ws.Range("U5").Value = 0.01807620809497196
$ws.Range("V5").Value = 0.05210415979209271
$ws.Range("W5").Value = 0
$ws.Range("X5").Value = 0.03619310018930463
$ws.Range("Y5").Value = 0.05182380314549435
$ws.Range("Z5").Value = 0.06776258362361258
$ws.Range("AA5").Value = 0
$ws.Range("AB5").Value = 0
$ws.Range("AC5").Value = 0
$ws.Range("AD5").Value = 0
$ws.Range("AE5").Value = 0
$ws.Range("AF5").Value = 0
$ws.Range("AG5").Value = 0
$ws.Range("AH5").Value = 0
$ws.Range("AI5").Value = 0
$ws.Range("B6").Value = 0
$ws.Range("C6").Value = 0
$ws.Range("D6").Value = 0.2096836020663077
$ws.Range("E6").Value = 0.1932369044015249
$ws.Range("F6").Value = 0.139790363997493
$ws.Range("G6").Value = 0.01138585923281594
$ws.Range("H6").Value = 0.01073764945913154
$ws.Range("I6").Value = 0
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 0
$ws.Range("L6").Value = 0.001579918808869742
$ws.Range("M6").Value = 0.1046423633812465
$ws.Range("N6").Value = 0.01979917202906401
$ws.Range("O6").Value = 0
$ws.Range("P6").Value = 0.02760780323836777
$ws.Range("Q6").Value = 0
$ws.Range("R6").Value = 0.1101760536828766
$ws.Range("S6").Value = 0
$ws.Range("T6").Value = 0.04736813344210433
$ws.Range("U6").Value = 0
$ws.Range("V6").Value = 0.01000598031887805
$ws.Range("W6").Value = 0.02623500148633265
$ws.Range("X6").Value = 0
$ws.Range("Y6").Value = 0.08546069235004941
$ws.Range("Z6").Value = 0
$ws.Range("AA6").Value = 0
$ws.Range("AB6").Value = 0
$ws.Range("AC6").Value = 0
$ws.Range("AD6").Value = 0.002290502104937686
$ws.Range("AE6").Value = 0
$ws.Range("AF6").Value = 0
$ws.Range("AG6").Value = 0
$ws.Range("AH6").Value = 0
$ws.Range("AI6").Value = 0

# --- Step2_Sj ---
$ws = $wb.Worksheets.Item("Step2_Sj")
$ws.Range("B2").Value = 0
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = 0.1608144049558392
$ws.Range("E2").Value = 0.1608144049558392
$ws.Range("F2").Value = 0.5311871414239212
$ws.Range("G2").Value = 0.5311871414239212
$ws.Range("H2").Value = 0.5410060973777523
$ws.Range("I2").Value = 0.5410060973777523
$ws.Range("J2").Value = 0.5410060973777523
$ws.Range("K2").Value = 0.5410060973777523
$ws.Range("L2").Value = 0.5839398225089902
$ws.Range("M2").Value = 0.5907153030175019
$ws.Range("N2").Value = 0.7031756575546299
$ws.Range("O2").Value = 0.7031756575546299
$ws.Range("P2").Value = 0.7141980808316174
$ws.Range("Q2").Value = 0.7141980808316174
$ws.Range("R2").Value = 0.7771908677625746
$ws.Range("S2").Value = 0.8331319312829895
$ws.Range("T2").Value = 0.8402053567038393
$ws.Range("U2").Value = 0.8402053567038393
$ws.Range("V2").Value = 0.884027163504846
$ws.Range("W2").Value = 0.884027163504846
$ws.Range("X2").Value = 0.8931932280189662
$ws.Range("Y2").Value = 0.9687863543041283
$ws.Range("Z2").Value = 1
$ws.Range("AA2").Value = 1
$ws.Range("AB2").Value = 1
$ws.Range("AC2").Value = 1
$ws.Range("AD2").Value = 1
$ws.Range("AE2").Value = 1
$ws.Range("AF2").Value = 1
$ws.Range("AG2").Value = 1
$ws.Range("AH2").Value = 1
$ws.Range("AI2").Value = 1
$ws.Range("B3").Value = 0
$ws.Range("C3").Value = 0
$ws.Range("D3").Value = 0.09328758731335501
$ws.Range("E3").Value = 0.09328758731335501
$ws.Range("F3").Value = 0.4399920219701463
$ws.Range("G3").Value = 0.4399920219701463
$ws.Range("H3").Value = 0.4399920219701463
$ws.Range("I3").Value = 0.4399920219701463
$ws.Range("J3").Value = 0.4399920219701463
$ws.Range("K3").Value = 0.4399920219701463
$ws.Range("L3").Value = 0.4831115566107278
$ws.Range("M3").Value = 0.4831115566107278
$ws.Range("N3").Value = 0.5955860736464588
$ws.Range("O3").Value = 0.614767933325998
$ws.Range("P3").Value = 0.614767933325998
$ws.Range("Q3").Value = 0.614767933325998
$ws.Range("R3").Value = 0.6305040873521965
$ws.Range("S3").Value = 0.7130686077973556
$ws.Range("T3").Value = 0.7130686077973556
$ws.Range("U3").Value = 0.7401297546843543
$ws.Range("V3").Value = 0.7742234805537999
$ws.Range("W3").Value = 0.7742234805537999
$ws.Range("X3").Value = 0.8383487688938911
$ws.Range("Y3").Value = 0.8605585325085947
$ws.Range("Z3").Value = 0.985517445098177
$ws.Range("AA3").Value = 0.985517445098177
$ws.Range("AB3").Value = 0.996971395659542
$ws.Range("AC3").Value = 0.996971395659542
$ws.Range("AD3").Value = 0.996971395659542
$ws.Range("AE3").Value = 0.9999999999999999
$ws.Range("AF3").Value = 0.9999999999999999
$ws.Range("AG3").Value = 0.9999999999999999
$ws.Range("AH3").Value = 0.9999999999999999
$ws.Range("AI3").Value = 0.9999999999999999
$ws.Range("B4").Value = 0
$ws.Range("C4").Value = 0
$ws.Range("D4").Value = 0.2229819650632389
$ws.Range("E4").Value = 0.2771222375436669
$ws.Range("F4").Value = 0.5114736435649256
$ws.Range("G4").Value = 0.5114736435649256
$ws.Range("H4").Value = 0.5114736435649256
$ws.Range("I4").Value = 0.5114736435649256
$ws.Range("J4").Value = 0.5114736435649256
$ws.Range("K4").Value = 0.5114736435649256
$ws.Range("L4").Value = 0.5604471421129897
$ws.Range("M4").Value = 0.5772548260331994
$ws.Range("N4").Value = 0.6733748623900401
$ws.Range("O4").Value = 0.6733748623900401
$ws.Range("P4").Value = 0.6772537214789753
$ws.Range("Q4").Value = 0.6772537214789753
$ws.Range("R4").Value = 0.7397306395388724
$ws.Range("S4").Value = 0.7732243481538832
$ws.Range("T4").Value = 0.7930003284883459
$ws.Range("U4").Value = 0.7930003284883459
$ws.Range("V4").Value = 0.8507418797707043
$ws.Range("W4").Value = 0.8573681672291446
$ws.Range("X4").Value = 0.8661881195031638
$ws.Range("Y4").Value = 0.9727261867800165
$ws.Range("Z4").Value = 0.9938593912302833
$ws.Range("AA4").Value = 0.9999999999999999
$ws.Range("AB4").Value = 0.9999999999999999
$ws.Range("AC4").Value = 0.9999999999999999
$ws.Range("AD4").Value = 0.9999999999999999
$ws.Range("AE4").Value = 0.9999999999999999
$ws.Range("AF4").Value = 0.9999999999999999
$ws.Range("AG4").Value = 0.9999999999999999
$ws.Range("AH4").Value = 0.9999999999999999
$ws.Range("AI4").Value = 0.9999999999999999
$ws.Range("B5").Value = 0
$ws.Range("C5").Value = 0
$ws.Range("D5").Value = 0.1235929690994866
$ws.Range("E5").Value = 0.1235929690994866
$ws.Range("F5").Value = 0.4759770868459081
$ws.Range("G5").Value = 0.4759770868459081
$ws.Range("H5").Value = 0.4869590415727605
$ws.Range("I5").Value = 0.4869590415727605
$ws.Range("J5").Value = 0.4869590415727605
$ws.Range("K5").Value = 0.4869590415727605
$ws.Range("L5").Value = 0.5222690242189132
$ws.Range("M5").Value = 0.5222690242189132
$ws.Range("N5").Value = 0.6539075925634277
$ws.Range("O5").Value = 0.6539075925634277
$ws.Range("P5").Value = 0.6539075925634277
$ws.Range("Q5").Value = 0.6599526755516272
$ws.Range("R5").Value = 0.6858617938169562
$ws.Range("S5").Value = 0.774040145154524
$ws.Range("T5").Value = 0.774040145154524
$ws.Range("U5").Value = 0.792116353249496
$ws.Range("V5").Value = 0.8442205130415887
$ws.Range("W5").Value = 0.8442205130415887
$ws.Range("X5").Value = 0.8804136132308934
$ws.Range("Y5").Value = 0.9322374163763878
$ws.Range("Z5").Value = 1
$ws.Range("AA5").Value = 1
$ws.Range("AB5").Value = 1
$ws.Range("AC5").Value = 1
$ws.Range("AD5").Value = 1
$ws.Range("AE5").Value = 1
$ws.Range("AF5").Value = 1
$ws.Range("AG5").Value = 1
$ws.Range("AH5").Value = 1
$ws.Range("AI5").Value = 1
$ws.Range("B6").Value = 0
$ws.Range("C6").Value = 0
$ws.Range("D6").Value = 0.2096836020663077
$ws.Range("E6").Value = 0.4029205064678326
$ws.Range("F6").Value = 0.5427108704653255
$ws.Range("G6").Value = 0.5540967296981415
$ws.Range("H6").Value = 0.5648343791572731
$ws.Range("I6").Value = 0.5648343791572731
$ws.Range("J6").Value = 0.5648343791572731
$ws.Range("K6").Value = 0.5648343791572731
$ws.Range("L6").Value = 0.5664142979661428
$ws.Range("M6").Value = 0.6710566613473894
$ws.Range("N6").Value = 0.6908558333764534
$ws.Range("O6").Value = 0.6908558333764534
$ws.Range("P6").Value = 0.7184636366148212
$ws.Range("Q6").Value = 0.7184636366148212
$ws.Range("R6").Value = 0.8286396902976977
$ws.Range("S6").Value = 0.8286396902976977
$ws.Range("T6").Value = 0.876007823739802
$ws.Range("U6").Value = 0.876007823739802
$ws.Range("V6").Value = 0.88601380405868
$ws.Range("W6").Value = 0.9122488055450126
$ws.Range("X6").Value = 0.9122488055450126
$ws.Range("Y6").Value = 0.997709497895062
$ws.Range("Z6").Value = 0.997709497895062
$ws.Range("AA6").Value = 0.997709497895062
$ws.Range("AB6").Value = 0.997709497895062
$ws.Range("AC6").Value = 0.997709497895062
$ws.Range("AD6").Value = 0.9999999999999997
$ws.Range("AE6").Value = 0.9999999999999997
$ws.Range("AF6").Value = 0.9999999999999997
$ws.Range("AG6").Value = 0.9999999999999997
$ws.Range("AH6").Value = 0.9999999999999997
$ws.Range("AI6").Value = 0.9999999999999997

# --- Step3_DataPts_0.5 ---
$ws = $wb.Worksheets.Item("Step3_DataPts_0.5")
$ws.Range("D2").Value = 5
$ws.Range("F2").Value = 0.5311871414239212
$ws.Range("G2").Value = 4
$ws.Range("K2").Value = '710R'
$ws.Range("D3").Value = 13
$ws.Range("F3").Value = 0.5955860736464588
$ws.Range("G3").Value = 12
$ws.Range("K3").Value = '710R'
$ws.Range("D4").Value = 5
$ws.Range("F4").Value = 0.5114736435649256
$ws.Range("G4").Value = 4
$ws.Range("K4").Value = '710R'
$ws.Range("D5").Value = 11
$ws.Range("F5").Value = 0.5222690242189132
$ws.Range("G5").Value = 10
$ws.Range("K5").Value = '710R'
$ws.Range("D6").Value = 5
$ws.Range("F6").Value = 0.5427108704653255
$ws.Range("G6").Value = 4
$ws.Range("K6").Value = '710R'

# --- Step3_DataPts_0.7 ---
$ws = $wb.Worksheets.Item("Step3_DataPts_0.7")
$ws.Range("D2").Value = 13
$ws.Range("F2").Value = 0.7031756575546299
$ws.Range("G2").Value = 12
$ws.Range("K2").Value = '710R'
$ws.Range("D3").Value = 18
$ws.Range("F3").Value = 0.7130686077973556
$ws.Range("G3").Value = 17
$ws.Range("K3").Value = '710R'
$ws.Range("D4").Value = 17
$ws.Range("F4").Value = 0.7397306395388724
$ws.Range("G4").Value = 16
$ws.Range("K4").Value = '710R'
$ws.Range("D5").Value = 18
$ws.Range("F5").Value = 0.774040145154524
$ws.Range("G5").Value = 17
$ws.Range("K5").Value = '710R'
$ws.Range("D6").Value = 15
$ws.Range("F6").Value = 0.7184636366148212
$ws.Range("G6").Value = 14
$ws.Range("K6").Value = '710R'

# --- Step3_DataPts_0.8 ---
$ws = $wb.Worksheets.Item("Step3_DataPts_0.8")
$ws.Range("D2").Value = 18
$ws.Range("F2").Value = 0.8331319312829895
$ws.Range("G2").Value = 17
$ws.Range("K2").Value = '710R'
$ws.Range("D3").Value = 23
$ws.Range("F3").Value = 0.8383487688938911
$ws.Range("G3").Value = 22
$ws.Range("K3").Value = '710R'
$ws.Range("D4").Value = 21
$ws.Range("F4").Value = 0.8507418797707043
$ws.Range("G4").Value = 20
$ws.Range("K4").Value = '710R'
$ws.Range("D5").Value = 21
$ws.Range("F5").Value = 0.8442205130415887
$ws.Range("G5").Value = 20
$ws.Range("K5").Value = '710R'
$ws.Range("D6").Value = 17
$ws.Range("F6").Value = 0.8286396902976977
$ws.Range("G6").Value = 16
$ws.Range("K6").Value = '710R'

# --- Step3_DataPts_0.9 ---
$ws = $wb.Worksheets.Item("Step3_DataPts_0.9")
$ws.Range("D2").Value = 24
$ws.Range("F2").Value = 0.9687863543041283
$ws.Range("G2").Value = 23
$ws.Range("K2").Value = '710R'
$ws.Range("D3").Value = 25
$ws.Range("F3").Value = 0.985517445098177
$ws.Range("G3").Value = 24
$ws.Range("K3").Value = '710R'
$ws.Range("D4").Value = 24
$ws.Range("F4").Value = 0.9727261867800165
$ws.Range("G4").Value = 23
$ws.Range("K4").Value = '710R'
$ws.Range("D5").Value = 24
$ws.Range("F5").Value = 0.9322374163763878
$ws.Range("G5").Value = 23
$ws.Range("K5").Value = '710R'
$ws.Range("D6").Value = 22
$ws.Range("F6").Value = 0.9122488055450126
$ws.Range("G6").Value = 21
$ws.Range("K6").Value = '710R'

Write-Host "Edit complete"